# Auto-generated Excel COM-interop script to apply the diff
# "Update latest output (run 202)"

$wb = $excel.ActiveWorkbook
$wsSchedule = $wb.Worksheets.Item("Schedule")
$wsDetailed = $wb.Worksheets.Item("Detailed")

# --- Sheet "Schedule": update row 2 and row 3 values, then remove rows 4 and 5 ---

$wsSchedule.Range("B2").Value = 46080.66666666666
$wsSchedule.Range("C2").Value = 16
$wsSchedule.Range("D2").Value = 60.48
$wsSchedule.Range("E2").Value = 1109.702724
$wsSchedule.Range("F2").Value = 18.34825932539682

$wsSchedule.Range("A3").Value = 46081.29166666666
$wsSchedule.Range("B3").Value = 46081.79166666666
$wsSchedule.Range("C3").Value = 12
$wsSchedule.Range("D3").Value = 45.36
$wsSchedule.Range("E3").Value = 300.791088
$wsSchedule.Range("F3").Value = 6.631196825396826

# Remove the now-obsolete rows 4 and 5 (dimension shrinks from A1:F5 to A1:F3)
$wsSchedule.Rows.Item(5).Delete()
$wsSchedule.Rows.Item(4).Delete()

# --- Sheet "Detailed": update individual cell values ---

$wsDetailed.Range("E12").Value = "ON"
$wsDetailed.Range("E13").Value = "ON"
$wsDetailed.Range("E14").Value = "ON"
$wsDetailed.Range("E15").Value = "ON"
$wsDetailed.Range("E16").Value = "ON"
$wsDetailed.Range("B38").Value = 50.46801
$wsDetailed.Range("B39").Value = 64.89
$wsDetailed.Range("B40").Value = 75.34523
$wsDetailed.Range("C40").Value = "historical"
$wsDetailed.Range("B41").Value = 72.90706
$wsDetailed.Range("C41").Value = "historical"
$wsDetailed.Range("B42").Value = 74.22528
$wsDetailed.Range("C42").Value = "historical"
$wsDetailed.Range("B43").Value = 78
$wsDetailed.Range("C43").Value = "historical"
$wsDetailed.Range("B44").Value = 71.40000000000001
$wsDetailed.Range("C44").Value = "historical"
$wsDetailed.Range("B45").Value = 65
$wsDetailed.Range("C45").Value = "historical"
$wsDetailed.Range("E45").Value = "OFF"
$wsDetailed.Range("B46").Value = 64.89
$wsDetailed.Range("C46").Value = "historical"
$wsDetailed.Range("E46").Value = "OFF"
$wsDetailed.Range("B47").Value = 65
$wsDetailed.Range("C47").Value = "historical"
$wsDetailed.Range("E47").Value = "OFF"
$wsDetailed.Range("B48").Value = 64.10590000000001
$wsDetailed.Range("C48").Value = "historical"
$wsDetailed.Range("E48").Value = "OFF"
$wsDetailed.Range("B49").Value = 65
$wsDetailed.Range("E49").Value = "OFF"
$wsDetailed.Range("B50").Value = 57.06
$wsDetailed.Range("E50").Value = "OFF"
$wsDetailed.Range("B51").Value = 65
$wsDetailed.Range("E51").Value = "OFF"
$wsDetailed.Range("B52").Value = 57.06
$wsDetailed.Range("E52").Value = "OFF"
$wsDetailed.Range("B53").Value = 57.06
$wsDetailed.Range("B54").Value = 56.98
$wsDetailed.Range("B55").Value = 56.15292
$wsDetailed.Range("B56").Value = 55.6532
$wsDetailed.Range("B57").Value = 54.99855
$wsDetailed.Range("B58").Value = 55.55142
$wsDetailed.Range("B59").Value = 55.8507
$wsDetailed.Range("B60").Value = 56.18178
$wsDetailed.Range("B61").Value = 57.31
$wsDetailed.Range("B62").Value = 57.36
$wsDetailed.Range("B63").Value = 57.36
$wsDetailed.Range("E64").Value = "ON"
$wsDetailed.Range("B66").Value = 1.17886
$wsDetailed.Range("B67").Value = 2.8337
$wsDetailed.Range("B68").Value = 1.53977
$wsDetailed.Range("B69").Value = 1.12995
$wsDetailed.Range("B70").Value = 6.75606
$wsDetailed.Range("B71").Value = 1.50663
$wsDetailed.Range("B72").Value = 0.7
$wsDetailed.Range("B73").Value = 1.46401
$wsDetailed.Range("B74").Value = 1.09962
$wsDetailed.Range("B75").Value = 1.45297
$wsDetailed.Range("B76").Value = 0.7
$wsDetailed.Range("B77").Value = 0.7
$wsDetailed.Range("B78").Value = 0.7
$wsDetailed.Range("B79").Value = 8.012169999999999
$wsDetailed.Range("B80").Value = 8.223660000000001
$wsDetailed.Range("B81").Value = 11.89595
$wsDetailed.Range("B82").Value = -3.76
$wsDetailed.Range("B83").Value = -4.14527
$wsDetailed.Range("B84").Value = 27.67105
$wsDetailed.Range("B86").Value = 50.73455
$wsDetailed.Range("E86").Value = "ON"
$wsDetailed.Range("B87").Value = 57.36
$wsDetailed.Range("E87").Value = "ON"
$wsDetailed.Range("B88").Value = 58.88215
$wsDetailed.Range("B89").Value = 63.60493
$wsDetailed.Range("B90").Value = 60.73634
$wsDetailed.Range("B91").Value = 58.92311
$wsDetailed.Range("B92").Value = 57.36
$wsDetailed.Range("B93").Value = 57.36
$wsDetailed.Range("B94").Value = 57.32
$wsDetailed.Range("B95").Value = 60.72274
$wsDetailed.Range("B97").Value = 57.06
